# Generate Report for Handback
# Update "Latest HO Xliff Generate Date" / handoff / handback timestamps
# for the c3bafcd9-... file row (row 3) across the Overview, zh-cn and
# de-de sheets, reflecting a newly generated handback report.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: row 3 -> c3bafcd9-f7bc-4b29-8611-a4694c15959a.md
# Column G = "Latest HO Xliff Generate Date"
$wsOverview.Range("G3").Value = "2016-10-18 04:13:42"

# zh-cn sheet: row 3 -> c3bafcd9-f7bc-4b29-8611-a4694c15959a.md
# Column H = "Correspond Handoff Datetime", Column K = "Correspond Handback DateTime"
$wsZhCn.Range("H3").Value = "2016-10-18 04:13:26"
$wsZhCn.Range("K3").Value = "2016-10-18 04:14:32"

# de-de sheet: row 3 -> c3bafcd9-f7bc-4b29-8611-a4694c15959a.md
# Column H = "Correspond Handoff Datetime", Column K = "Correspond Handback DateTime"
$wsDeDe.Range("H3").Value = "2016-10-18 04:13:42"
$wsDeDe.Range("K3").Value = "2016-10-18 04:14:58"
